$d = $word.ActiveDocument

# Tech Stack lines: each occurrence originally reads
#   ... <run-with-colon>:</run> <run2>" Python, Apache Airflow, AWS, MSSQL, Linux, Shell Scripting, Internal VPS, "</run2> <run3>"CI/CD"</run3>
# and needs to become
#   ... <run-with-colon>:</run> <run2>" "</run2> <run3>"Python, Apache Airflow (Workflow Orchestration / CI/CD), AWS, MSSQL, Linux, Bash, On-prem VPS"</run3>
# (project 1), and similar replacements for projects 2 and 3 below.
#
# We replace run2's and run3's text in place (not touching the preceding
# ":"/"Stack:" run) so the paragraph keeps the same 3-run shape as the
# target instead of Word's Find/Replace collapsing everything into one
# run. A Bold-toggle bracket around each write keeps the edited run from
# being silently re-merged with its neighbours (which otherwise happens
# whenever two adjacent runs end up with identical formatting).

$run2Old = " Python, Apache Airflow, AWS, MSSQL, Linux, Shell Scripting, Internal VPS, "
$run3Old = "CI/CD"

$replacements = @(
    "Python, Apache Airflow (Workflow Orchestration / CI/CD), AWS, MSSQL, Linux, Bash, On-prem VPS",
    "Python, Custom Distributed Processing (VMware-based VPS Cluster), MSSQL, Pandas, NumPy",
    "Python, pymssql, Custom Distributed Processing, Multiprocessing, XQuery"
)

$searchStart = 0
foreach ($run3New in $replacements) {
    $run2New = " "

    $rng = $d.Range($searchStart, $d.Content.End)
    $found = $rng.Find.Execute($run2Old + $run3Old, $true, $false, $false, $false, $false, $true, 1, $false)
    if (-not $found) {
        Write-Host "NOT FOUND starting at $searchStart"
        continue
    }

    $S = $rng.Start
    $E = $rng.End
    $splitPoint = $S + $run2Old.Length

    # --- run3 ("CI/CD" -> full new stack text) ---
    $r3 = $d.Range($splitPoint, $E)
    $r3.Bold = 1
    $r3.Text = $run3New
    $r3b = $d.Range($splitPoint, $splitPoint + $run3New.Length)
    $r3b.Bold = 0

    # --- run2 (the long " Python, ... Internal VPS, " run -> single space) ---
    $r2 = $d.Range($S, $splitPoint)
    $r2.Bold = 1
    $r2.Text = $run2New
    $r2b = $d.Range($S, $S + $run2New.Length)
    $r2b.Bold = 0

    $searchStart = $S + $run2New.Length + $run3New.Length
    Write-Host "Replaced occurrence at $S"
}

Write-Host "done"
